# Insert a new data row at row 3 (pushing existing rows 3..83 down to 4..84)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(3).Insert()

# Populate the newly inserted row 3 with its data
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C3").Value = "Arica y Parinacota"
$ws.Range("D3").Value = 44812
$ws.Range("E3").Value = 15
$ws.Range("F3").Value = 100114001
$ws.Range("G3").Value = "Papa"
$ws.Range("H3").Value = "Rodeo"
$ws.Range("I3").Value = "1a (guarda)"
$ws.Range("J3").Value = 1000
$ws.Range("K3").Value = 8000
$ws.Range("L3").Value = 9000
$ws.Range("M3").Value = 8500
$ws.Range("N3").Value = "`$/saco 25 kilos"
$ws.Range("O3").Value = "Región de Los Lagos"
$ws.Range("P3").Value = 340
$ws.Range("Q3").Value = 25
$ws.Range("R3").Value = "Hortaliza"
